# Refresh the LR-pair sheet with newly-computed TPM-based values.
# - Two self-pairing rows (sending cluster == target cluster == MuSCs,
#   and the old ECs-as-sender / MuSCs-as-target row) are dropped, so the
#   sheet shrinks from 6 data rows (rows 2-7) to 4 data rows (rows 2-5).
# - Every remaining numeric column (E:T) is recomputed with new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that no longer exist in the refreshed data.
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows("7:7").Delete()
$ws.Rows("4:4").Delete()

# After the deletes the sheet holds, in order:
#   row2 = ECs   -> ECs
#   row3 = ECs   -> FAPs
#   row4 = MuSCs -> ECs   (was row 5 before the delete)
#   row5 = MuSCs -> FAPs  (was row 6 before the delete)
# Labels (columns A-D) are already correct; only E:T need new values.

$rowValues = @{
    2 = @(1, 0.3333333333333333, 0.009651333333333333, 0.028954, 0.8567539576860482, 0.8567539576860482, 1, 0.3333333333333333, 0.006768333333333334, 0.020305, 0.00347174015482542, 0.00347174015482542, 0.00006532344111111111, 0.00058791097, 0.002974427117704252, 0.002974427117704252)
    3 = @(1, 0.3333333333333333, 0.009651333333333333, 0.028954, 0.8567539576860482, 0.8567539576860482, 3, 1, 1.942782333333333, 5.828347, 0.9965282598451746, 0.9965282598451747, 0.01875043989311111, 0.168753959038, 0.853779530568344, 0.8537795305683441)
    4 = @(1, 0.3333333333333333, 0.001613666666666667, 0.004841, 0.1432460423139518, 0.1432460423139518, 1, 0.3333333333333333, 0.006768333333333334, 0.020305, 0.00347174015482542, 0.00347174015482542, 0.00001092183388888889, 0.000098296505, 0.0004973130371211676, 0.0004973130371211676)
    5 = @(1, 0.3333333333333333, 0.001613666666666667, 0.004841, 0.1432460423139518, 0.1432460423139518, 3, 1, 1.942782333333333, 5.828347, 0.9965282598451746, 0.9965282598451747, 0.003135003091888889, 0.028215027827, 0.1427487292768306, 0.1427487292768306)
}

foreach ($r in $rowValues.Keys) {
    $vals = $rowValues[$r]
    $col = 5
    foreach ($v in $vals) {
        $ws.Cells.Item($r, $col).Value = $v
        $col = $col + 1
    }
}
